$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 21 (2025Q3) metrics per updated "bibi e add" data
$ws.Range("C21").Value = 269
$ws.Range("D21").Value = 236
$ws.Range("E21").Value = 33
$ws.Range("F21").Value = 67.621776504298
